$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dates stays text, B2:F2 updated numeric values, G2/H2 unchanged
$ws.Range("A2").Value = "30/05/2000"
$ws.Range("B2").Value = 4000
$ws.Range("C2").Value = 4000
$ws.Range("D2").Value = 4000
$ws.Range("E2").Value = 4000
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 100

# Row 3: date text, and B3:H3 become text-typed numeric-looking strings
$ws.Range("A3").Value = "31/07/2000"

$ws.Range("B3:H3").NumberFormat = "@"
$ws.Range("B3").Value = "1000.00"
$ws.Range("C3").Value = "5000.00"
$ws.Range("D3").Value = "1000.00"
$ws.Range("E3").Value = "5000.00"
$ws.Range("F3").Value = "60.0"
$ws.Range("G3").Value = "0.00"
$ws.Range("H3").Value = "100.00"
$ws.Range("B3:H3").ClearFormats()

# Remove rows 4-6 (they no longer exist in the updated table)
$ws.Range("A4:H6").EntireRow.Delete()
